# Add 2022-Q4 data:
#  1. Insert a new row (row 2) into the "总计" (total) summary sheet for 2022-Q4.
#  2. Insert a brand-new "2022-Q4" worksheet (positioned right after "总计",
#     before "2022-Q3") with the quarter's fund-holdings detail table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert new row 2 for 2022-Q4 (shifts existing rows down).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 16
$total.Cells.Item(2, 4).Value = 6.97

# Re-apply the same look as the other index cells in column A (bold, boxed,
# centered) by copying the format from the row below, which still carries it.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A is a plain sequential row index (0,1,2,...), independent of the
# quarter label - renumber the rows that got pushed down by the insert.
for ($i = 3; $i -le 8; $i++) {
    $total.Cells.Item($i, 1).Value = $i - 2
}

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund table, inserted before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)

# Clone the existing "2022-Q3" sheet (keeps header/column styling identical)
# and place the clone immediately before it, then rename it.
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The source sheet only has 4 data rows (rows 2-5); the 2022-Q4 table needs
# 16, so extend the styled block (row 5's formatting) down through row 17.
$q4.Rows.Item(5).Copy()
$q4.Range("A6:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @(0,  "163402", "兴全趋势投资混合（LOF）",        "205.32", "87.99", "2.86", "5.8722", 10),
    @(1,  "001445", "华安国企改革主题灵活配置混合A",   "15.25",  "85.79", "2.61", "0.3980", 10),
    @(2,  "006058", "民生加银新兴成长混合",            "3.44",   "86.91", "3.83", "0.1318", 9),
    @(3,  "013676", "兴银兴慧一年持有混合A",           "8.13",   "23.86", "1.16", "0.0943", 3),
    @(4,  "501200", "民生加银科技创新 3 年封闭混合",    "2.45",   "86.66", "3.67", "0.0899", 10),
    @(5,  "010122", "华泰柏瑞优势领航混合A",           "3.69",   "93.67", "2.37", "0.0875", 7),
    @(6,  "013677", "兴银兴慧一年持有混合C",           "4.59",   "23.86", "1.16", "0.0532", 3),
    @(7,  "519644", "银河智联主题灵活配置混合",         "1.20",   "89.42", "4.38", "0.0526", 10),
    @(8,  "009206", "兴银丰运稳益回报混合C",           "3.03",   "39.08", "1.67", "0.0506", 4),
    @(9,  "005041", "人保研究精选混合A",               "1.09",   "87.60", "2.98", "0.0325", 8),
    @(10, "014839", "兴银碳中和主题混合C",             "0.64",   "92.17", "5.07", "0.0324", 3),
    @(11, "009205", "兴银丰运稳益回报混合A",           "1.91",   "39.08", "1.67", "0.0319", 4),
    @(12, "014838", "兴银碳中和主题混合A",             "0.53",   "92.17", "5.07", "0.0269", 3),
    @(13, "010123", "华泰柏瑞优势领航混合C",           "0.37",   "93.67", "2.37", "0.0088", 7),
    @(14, "016290", "华安国企改革主题灵活配置混合C",    "0.16",   "85.79", "2.61", "0.0042", 10),
    @(15, "005042", "人保研究精选混合C",               "0.05",   "87.60", "2.98", "0.0015", 8)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $q4.Range("B" + $r + ":G" + $r).ClearFormats()
    $r = $r + 1
}
